# Update WBS/requirements workbook: add requirement references (column F)
# for each task row, and refresh the column widths / selection state to
# match the latest authoring session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (col B, C, D, F) -----------------------------------
# ColumnWidth is specified in "characters"; the underlying XML width is
# quantized in 1/6-character increments by the engine, so we back-solve
# the input that lands closest to the desired stored width.
$ws.Columns.Item(2).ColumnWidth = 23.666666666666668   # -> stored width 24.5
$ws.Columns.Item(3).ColumnWidth = 30.830729166666668   # -> stored width ~31.6667 (closest to 31.6640625)
$ws.Columns.Item(4).ColumnWidth = 25.666666666666668   # -> stored width 26.5
$ws.Columns.Item(6).ColumnWidth = 56.666666666666664   # -> stored width 57.5

# --- New "REQUISITO" references in column F ---------------------------
# Values are entered in the order that reproduces the original shared
# string table layout: REQ_1/2, REQ_4, REQ_6, REQ_3/5, REQ_5, REQ_8,
# REQ_7, REQ_9, REQ_11, REQ_10.
$ws.Range("F11").Value = "REQ_1, REQ_2"
$ws.Range("F23").Value = "REQ_4"
$ws.Range("F26").Value = "REQ_6"
$ws.Range("F13").Value = "REQ_3, REQ_5"
$ws.Range("F15").Value = "REQ_5"
$ws.Range("F16").Value = "REQ_8"
$ws.Range("F17").Value = "REQ_7"
$ws.Range("F12").Value = "REQ_9"
$ws.Range("F19").Value = "REQ_11"
$ws.Range("F14").Value = "REQ_10"
$ws.Range("F24").Value = "REQ_4"
$ws.Range("F20").Value = "REQ_11"
$ws.Range("F21").Value = "REQ_11"

# --- Refresh the active selection on the sheet -------------------------
$ws.Range("F14").Select() | Out-Null
